# Timesheet update: period May -> June, swap activity descriptions for rows 12/13,
# and fill in the daily hour tallies for rows 11-13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text / label updates -------------------------------------------------

# Row 12 & 13 activity descriptions (row 11 text is unchanged)
$ws.Range("B12").Value = "Support SIT Release 8 (Billing)"
$ws.Range("B13").Value = "Support  UAT Release 8 (Billing)"

# Period label (month/year header area)
$ws.Range("O6").Value = "20 May - 19 June"

# --- Row 11: "Diskusi dengan tim developer Billing" hour tally -----------

$ws.Range("C11").Value = 1
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 1
$ws.Range("J11").Value = 1
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 1
$ws.Range("N11").Value = 1
$ws.Range("Q11").Value = 1
$ws.Range("R11").Value = 1
$ws.Range("S11").Value = 1
$ws.Range("T11").Value = 1
$ws.Range("U11").Value = 1
$ws.Range("X11").Value = 1
$ws.Range("Y11").Value = 1
$ws.Range("Z11").Value = 1
$ws.Range("AA11").Value = 1
$ws.Range("AB11").Value = 1
$ws.Range("AG11").Value = 1

# --- Row 12: "Support SIT Release 8 (Billing)" hour tally ----------------

$ws.Range("C12").Value = 6
$ws.Range("D12").Value = 6
$ws.Range("E12").Value = 6
$ws.Range("J12").Value = 6
$ws.Range("K12").Value = 6
$ws.Range("L12").Value = 6
$ws.Range("M12").Value = 6
$ws.Range("N12").Value = 6
$ws.Range("Q12").Value = 6
$ws.Range("R12").Value = 6
$ws.Range("S12").Value = 6
$ws.Range("T12").Value = 6
$ws.Range("U12").Value = 6
$ws.Range("X12").Value = 5
$ws.Range("Y12").Value = 5
$ws.Range("Z12").Value = 5
$ws.Range("AA12").Value = 5
$ws.Range("AB12").Value = 5

# --- Row 13: "Support  UAT Release 8 (Billing)" hour tally ---------------
# (previously held the Z13/AA13/AB13 = 5 values that now move/shrink)

$ws.Range("Z13").Value = 1
$ws.Range("AA13").Value = 1
$ws.Range("AB13").Value = 1
$ws.Range("X13").Value = 1
$ws.Range("Y13").Value = 1
$ws.Range("AG13").Value = 6

# --- View bookkeeping (matches the saved selection/scroll position) ------

$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("AH13").Select()
